$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C labels
$ws.Range("C1").Value = "average"
$ws.Range("C2").Value = "max"
$ws.Range("C3").Value = "min"

# Column D summary formulas over the 100 simulation results in column A
$ws.Range("D1").Formula = "=AVERAGE(A:A)"
$ws.Range("D2").Formula = "=MAX(A:A)"
$ws.Range("D3").Formula = "=MIN(A:A)"

# Highlight the average value in bold
$ws.Range("D1").Font.Bold = $true

# Match the saved selection state (C1:D3, active cell C1)
$ws.Range("C1:D3").Select() | Out-Null
